$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row of data: 06.12.2016, 13, 14, 7
# Force the date-like text to be stored as a literal string (matching how
# the existing date text cells are stored) rather than being auto-parsed
# into a date serial value.
$cell = $ws.Range("A5")
$cell.NumberFormat = "@"
$cell.Value = "06.12.2016"
$cell.Style = "Normal"

$ws.Range("B5").Value = 13
$ws.Range("C5").Value = 14
$ws.Range("D5").Value = 7

# Update selection to match the target state
$ws.Range("D7").Select()
